# Auto-generated edit script: updates cryptos worksheet values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''90.210.83'
$ws.Range("E2").Value = '  -0.76%  '

$ws.Range("D3").Value = '''3.083.69'
$ws.Range("E3").Value = '  -2.71%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").Value = '''232.93'
$ws.Range("E5").Value = '  +6.99%  '

$ws.Range("D6").Value = '''624.64'
$ws.Range("E6").Value = '  -0.53%  '

$ws.Range("E7").Value = '  -4.05%  '

$ws.Range("D8").Value = '''0.360'
$ws.Range("E8").Value = '  -2.64%  '

$ws.Range("E9").Value = '  +0.00%  '

$ws.Range("D10").Value = '''0.726'
$ws.Range("E10").Value = '  -4.80%  '

$ws.Range("D11").Value = '''2.498.04'
$ws.Range("E11").Value = '  -21.06%  '

$ws.Range("D12").Value = '''0.196'
$ws.Range("E12").Value = '  -3.07%  '

$ws.Range("D13").Value = '''36.29'
$ws.Range("E13").Value = '  +2.80%  '

$ws.Range("E14").Value = '  -0.22%  '

$ws.Range("D15").Value = '''5.46'
$ws.Range("E15").Value = '  -4.51%  '

$ws.Range("D16").Value = '''90.118.79'
$ws.Range("E16").Value = '  -0.65%  '

$ws.Range("D17").Value = '''3.658.96'
$ws.Range("E17").Value = '  -2.45%  '

$ws.Range("D18").Value = '''3.088.00'
$ws.Range("E18").Value = '  -2.44%  '

$ws.Range("E19").Value = '  +1.14%  '

$ws.Range("E20").Value = '  -2.68%  '

$ws.Range("D21").Value = '''14.02'
$ws.Range("E21").Value = '  -2.35%  '

$ws.Range("D22").Value = '''438.25'
$ws.Range("E22").Value = '  -2.74%  '

$ws.Range("D23").Value = '''5.56'
$ws.Range("E23").Value = '  +6.39%  '

$ws.Range("D24").Value = '''8.87'
$ws.Range("E24").Value = '  -1.56%  '

$ws.Range("D25").Value = '''7.55'
$ws.Range("E25").Value = '  -1.99%  '

$ws.Range("D26").Value = '''5.68'
$ws.Range("E26").Value = '  -4.92%  '

$ws.Range("D27").Value = '''89.10'
$ws.Range("E27").Value = '  -1.69%  '

$ws.Range("D28").Value = '''12.16'
$ws.Range("E28").Value = '  -0.28%  '

$ws.Range("D29").Value = '''3.291.19'
$ws.Range("E29").Value = '  -1.47%  '

$ws.Range("E30").Value = '  -0.02%  '

$ws.Range("E31").Value = '  +1.15%  '

$ws.Range("E32").Value = '  -2.48%  '

$ws.Range("D33").Value = '''0.976'
$ws.Range("E33").Value = '  -3.83%  '

$ws.Range("E34").Value = '  +12.48%  '

$ws.Range("D35").Value = '''26.25'
$ws.Range("E35").Value = '  +1.36%  '

$ws.Range("D36").Value = '''0.153'
$ws.Range("E36").Value = '  +5.70%  '

$ws.Range("D37").Value = '''3.79'
$ws.Range("E37").Value = '  +2.29%  '

$ws.Range("D38").Value = '''503.86'
$ws.Range("E38").Value = '  -4.28%  '

$ws.Range("E39").Value = '  -0.15%  '

$ws.Range("D40").Value = '''6.96'
$ws.Range("E40").Value = '  -0.43%  '

$ws.Range("E41").Value = '  -2.02%  '

$ws.Range("D42").Value = '''0.0895'
$ws.Range("E42").Value = '  +0.66%  '

$ws.Range("B43").Value = 'PolygonEcosystemToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D43").Value = '''0.410'
$ws.Range("E43").Value = '  -3.15%  '

$ws.Range("B44").Value = 'MantraDAO'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D44").Value = '''3.51'
$ws.Range("E44").Value = '  +52.66%  '

$ws.Range("D45").Value = '''22.17'
$ws.Range("E45").Value = '  -0.24%  '

$ws.Range("E47").Value = '  -2.78%  '

$ws.Range("D48").Value = '''150.52'
$ws.Range("E48").Value = '  +1.73%  '

$ws.Range("D49").Value = '''0.687'
$ws.Range("E49").Value = '  +4.40%  '

$ws.Range("D50").Value = '''45.01'
$ws.Range("E50").Value = '  +1.57%  '

$ws.Range("E51").Value = '  -0.62%  '

